# Applies the crypto price/volume refresh described in the commit.
# Column D values are numeric-looking strings (e.g. "58.722.08") that must
# stay as literal Text, so NumberFormat is forced to "@" before assignment
# (otherwise Excel silently re-parses them as numbers/dates and mangles them,
# e.g. "523.10" -> 523.1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.722.08"
$ws.Range("E2").Value = "  +0.89%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.624.79"
$ws.Range("E3").Value = "  +1.93%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "523.10"
$ws.Range("E5").Value = "  +3.40%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.23"
$ws.Range("E6").Value = "  -0.02%  "

# Row 7
$ws.Range("E7").Value = "  +0.11%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.568"
$ws.Range("E8").Value = "  -0.07%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.631.76"
$ws.Range("E9").Value = "  +1.42%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.27"
$ws.Range("E10").Value = "  -0.88%  "

# Row 11
$ws.Range("E11").Value = "  +0.65%  "

# Row 12
$ws.Range("E12").Value = "  -0.33%  "

# Row 13
$ws.Range("E13").Value = "  -0.97%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.087.07"
$ws.Range("E14").Value = "  +1.86%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.755.50"
$ws.Range("E15").Value = "  +0.98%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.72"
$ws.Range("E16").Value = "  -1.10%  "

# Row 17
$ws.Range("E17").Value = "  -0.61%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.634.87"
$ws.Range("E18").Value = "  +1.83%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "344.50"
$ws.Range("E19").Value = "  +0.75%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.41"
$ws.Range("E20").Value = "  -2.64%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.15"
$ws.Range("E21").Value = "  -1.12%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.11"
$ws.Range("E22").Value = "  +0.94%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.04%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.33"
$ws.Range("E24").Value = "  +1.24%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.413"
$ws.Range("E25").Value = "  -0.99%  "

# Row 26
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.38%  "

# Row 27
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.163"
$ws.Range("E27").Value = "  +3.18%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0795"
$ws.Range("E28").Value = "  -2.20%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.94"
$ws.Range("E29").Value = "  -0.61%  "

# Row 30
$ws.Range("E30").Value = "  +0.05%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.18"
$ws.Range("E31").Value = "  +1.21%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.58"
$ws.Range("E32").Value = "  +2.81%  "

# Row 33
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.81"
$ws.Range("E33").Value = "  +0.17%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.48"
$ws.Range("E34").Value = "  +1.65%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.970"
$ws.Range("E35").Value = "  +1.05%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.96"
$ws.Range("E36").Value = "  -0.18%  "

# Row 37
$ws.Range("E37").Value = "  +0.35%  "

# Row 38
$ws.Range("E38").Value = "  +1.82%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.833"
$ws.Range("E39").Value = "  -1.61%  "

# Row 40
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.62"
$ws.Range("E40").Value = "  +1.20%  "

# Row 41
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.41"
$ws.Range("E41").Value = "  +1.36%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "277.78"
$ws.Range("E42").Value = "  -3.04%  "

# Row 43
$ws.Range("E43").Value = "  +0.05%  "

# Row 44
$ws.Range("E44").Value = "  -0.86%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.598"
$ws.Range("E45").Value = "  -2.06%  "

# Row 46
$ws.Range("E46").Value = "  +1.80%  "

# Row 47
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.31"
$ws.Range("E47").Value = "  +0.59%  "

# Row 48
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0519"
$ws.Range("E48").Value = "  -2.75%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.986.55"
$ws.Range("E49").Value = "  +2.65%  "

# Row 50
$ws.Range("E50").Value = "  +0.35%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.62"
$ws.Range("E51").Value = "  +1.52%  "
